$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 4690.4
$ws.Range("I34").Value = 4690.4
$ws.Range("K34").Value = 4690.4
$ws.Range("M34").Value = -4487.4
$ws.Range("H36").Value = 4690.4
$ws.Range("I36").Value = 4690.4
$ws.Range("K36").Value = 4690.4
$ws.Range("M36").Value = -3975.4
$ws.Range("H40").Value = 4983.8
$ws.Range("J40").Value = 5308.2573
$ws.Range("L40").Value = 5308.2573
$ws.Range("N40").Value = -5658.2573
$ws.Range("H43").Value = 3267.3333
$ws.Range("J43").Value = 3651
$ws.Range("L43").Value = 3651
$ws.Range("N43").Value = -3789
$ws.Range("H47").Value = 3500
$ws.Range("I47").Value = 3500
$ws.Range("K47").Value = 3500
$ws.Range("M47").Value = -2528
$ws.Range("H54").Value = 6300
$ws.Range("J54").Value = 6300
$ws.Range("L54").Value = 6300
$ws.Range("N54").Value = -7272
$ws.Range("H76").Value = 4243.9414
$ws.Range("I76").Value = 3938.5
$ws.Range("K76").Value = 3938.5
$ws.Range("M76").Value = -3623.5
$ws.Range("H79").Value = 4243.9414
$ws.Range("I79").Value = 3938.5
$ws.Range("K79").Value = 3938.5
$ws.Range("M79").Value = -2846.5
$ws.Range("H98").Value = 2961
$ws.Range("J98").Value = 2388
$ws.Range("L98").Value = 2388
$ws.Range("N98").Value = -5384
$ws.Range("H106").Value = 5878.154
$ws.Range("I106").Value = 5784.6665
$ws.Range("K106").Value = 5784.6665
$ws.Range("M106").Value = -5153.6665
$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("M114").ClearContents()
$ws.Range("H122").Value = 2961
$ws.Range("J122").Value = 2388
$ws.Range("L122").Value = 7164
$ws.Range("N122").Value = -12064
$ws.Range("H126").Value = 69499.5
$ws.Range("J126").Value = 69499.5
$ws.Range("L126").Value = 69499.5
$ws.Range("N126").Value = -79379.5
$ws.Range("H132").Value = 8476.629999999999
$ws.Range("I132").Value = 8571.885
$ws.Range("K132").Value = 25715.655
$ws.Range("M132").Value = -23185.655
$ws.Range("H137").Value = 38474664
$ws.Range("I137").Value = 55557868
$ws.Range("K137").Value = 166673604
$ws.Range("M137").Value = -166671054

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1466422.6
$ws.Range("I32").Value = 1589108.8
$ws.Range("K32").Value = 1589108.8
$ws.Range("M32").Value = -1588821.8
$ws.Range("H61").Value = 5009793
$ws.Range("I61").Value = 15613.875
$ws.Range("J61").Value = 8339245.5
$ws.Range("K61").Value = 15613.875
$ws.Range("L61").Value = 8339245.5
$ws.Range("M61").Value = -15401.875
$ws.Range("N61").Value = -8339669.5
$ws.Range("H63").Value = 2391.8
$ws.Range("J63").Value = 2879.9
$ws.Range("L63").Value = 2879.9
$ws.Range("N63").Value = -4251.9
$ws.Range("H66").Value = 2391.8
$ws.Range("J66").Value = 2879.9
$ws.Range("L66").Value = 14399.5
$ws.Range("N66").Value = -21263.5
$ws.Range("H74").Value = 2657744.5
$ws.Range("I74").Value = 3276271.8
$ws.Range("K74").Value = 3276271.8
$ws.Range("M74").Value = -3275397.8
$ws.Range("H77").Value = 2657744.5
$ws.Range("I77").Value = 3276271.8
$ws.Range("K77").Value = 16381359
$ws.Range("M77").Value = -16376991
$ws.Range("H110").Value = 1518.7931
$ws.Range("I110").Value = 1143.4375
$ws.Range("K110").Value = 1143.4375
$ws.Range("M110").Value = 901.5625
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("L119").ClearContents()
$ws.Range("M119").ClearContents()
$ws.Range("N119").Value = 0
$ws.Range("H136").Value = 5009793
$ws.Range("I136").Value = 15613.875
$ws.Range("J136").Value = 8339245.5
$ws.Range("K136").Value = 46841.625
$ws.Range("L136").Value = 25017736.5
$ws.Range("M136").Value = -44291.625
$ws.Range("N136").Value = -25022836.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4566.875
$ws.Range("I94").Value = 4566.875
$ws.Range("K94").Value = 4566.875
$ws.Range("M94").Value = -4115.875
$ws.Range("H97").Value = 35476
$ws.Range("I97").Value = 23214
$ws.Range("K97").Value = 23214
$ws.Range("M97").Value = -22223
$ws.Range("H99").Value = 6980.6113
$ws.Range("J99").Value = 2499.6667
$ws.Range("L99").Value = 2499.6667
$ws.Range("N99").Value = -5495.6667
$ws.Range("H105").Value = 5277.381
$ws.Range("I105").Value = 2982
$ws.Range("K105").Value = 2982
$ws.Range("M105").Value = -1235

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 84.15000000000001
$ws.Range("I7").Value = 69.666664
$ws.Range("K7").Value = 69.666664
$ws.Range("M7").Value = 43.333336
$ws.Range("H16").Value = 227522.2
$ws.Range("I16").Value = 227522.2
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 227522.2
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -227235.2
$ws.Range("H31").Value = 5853145
$ws.Range("I31").Value = 6584538
$ws.Range("K31").Value = 6584538
$ws.Range("M31").Value = -6584243
$ws.Range("H34").Value = 5853145
$ws.Range("I34").Value = 6584538
$ws.Range("K34").Value = 6584538
$ws.Range("M34").Value = -6584336
$ws.Range("H113").Value = 227522.2
$ws.Range("I113").Value = 227522.2
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 227522.2
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -225352.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 524.8
$ws.Range("I12").Value = 1131.6
$ws.Range("J12").Value = 221.4
$ws.Range("K12").Value = 3394.8
$ws.Range("L12").Value = 664.2
$ws.Range("M12").Value = -3221.8
$ws.Range("N12").Value = -1010.2
$ws.Range("H56").Value = 9727
$ws.Range("I56").Value = 9727
$ws.Range("K56").Value = 9727
$ws.Range("M56").Value = -9197
$ws.Range("H126").Value = 15000
$ws.Range("J126").Value = 15000
$ws.Range("L126").Value = 45000
$ws.Range("N126").Value = -54880
$ws.Range("H131").Value = 5063.878
$ws.Range("I131").Value = 2935.25
$ws.Range("J131").Value = 5294
$ws.Range("K131").Value = 8805.75
$ws.Range("L131").Value = 15882
$ws.Range("M131").Value = -3765.75
$ws.Range("N131").Value = -25962
$ws.Range("H132").Value = 2360.5833
$ws.Range("I132").Value = 1993.3636
$ws.Range("J132").Value = 6400
$ws.Range("K132").Value = 17940.2724
$ws.Range("L132").Value = 57600
$ws.Range("M132").Value = -15410.2724
$ws.Range("N132").Value = -62660
$ws.Range("H140").Value = 2771.7273
$ws.Range("I140").Value = 2748.9
$ws.Range("K140").Value = 8246.700000000001
$ws.Range("M140").Value = -3066.700000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1599.5
$ws.Range("I80").Value = 1599
$ws.Range("K80").Value = 1599
$ws.Range("M80").Value = -601
$ws.Range("H83").Value = 1599.5
$ws.Range("I83").Value = 1599
$ws.Range("K83").Value = 7995
$ws.Range("M83").Value = -3003
$ws.Range("H97").Value = 1195.6
$ws.Range("I97").Value = 1166.0416
$ws.Range("K97").Value = 1166.0416
$ws.Range("M97").Value = -670.0416
$ws.Range("H102").Value = 1149.5
$ws.Range("I102").Value = 1149.5
$ws.Range("K102").Value = 1149.5
$ws.Range("M102").Value = 472.5
$ws.Range("H113").Value = 1772.7273
$ws.Range("J113").Value = 3000
$ws.Range("L113").Value = 3000
$ws.Range("N113").Value = -7340
$ws.Range("H128").Value = 70000
$ws.Range("J128").Value = 70000
$ws.Range("L128").Value = 70000
$ws.Range("N128").Value = -79960

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3446.2334
$ws.Range("J46").Value = 5028.647
$ws.Range("L46").Value = 5028.647
$ws.Range("N46").Value = -5404.647
$ws.Range("H93").Value = 9594.4
$ws.Range("I93").Value = 5251.5
$ws.Range("K93").Value = 5251.5
$ws.Range("M93").Value = -4003.5
$ws.Range("H122").Value = 3156.2307
$ws.Range("I122").Value = 2821
$ws.Range("K122").Value = 8463
$ws.Range("M122").Value = -6013
